$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.737949132919312
$ws.Range("B1").Value = 2.282087087631226
$ws.Range("C1").Value = 2.336277484893799
$ws.Range("D1").Value = 2.64025092124939
$ws.Range("E1").Value = 3.37243914604187
